# "updated single child outputs"
# The sheet had an extra (duplicate) last row (row 6, a repeat of the
# even_MAG-GUT71751.fa row) and placeholder prediction-probability values
# of 1 in column B. This normalises the sheet back down to the five real
# rows and fills column B with the real per-row predicted scores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the duplicated trailing row (old row 6: even_MAG-GUT71751.fa again).
$ws.Rows.Item(6).Delete()

# The surviving rows should all share the same (already-existing) text
# style used by the old row 6, instead of the separate duplicate style
# the first five rows used to carry.
$ws.Range("A1:A5").NumberFormat = "@"
$ws.Range("B1").NumberFormat = "@"
$ws.Range("C1").NumberFormat = "@"

# Replace the placeholder "1" values in column B with the real scores.
$ws.Range("B2").Value2 = -6.879024153281895
$ws.Range("B3").Value2 = -3.0349666117530116
$ws.Range("B4").Value2 = -1.8225386348122488
$ws.Range("B5").Value2 = -1.2184361197441227
